$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 14.80997966666667
$ws.Range("H2").Value = 44.429939
$ws.Range("I2").Value = 0.2388798507865045
$ws.Range("J2").Value = 0.2388798507865045
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 5.872120333333332
$ws.Range("N2").Value = 17.616361
$ws.Range("O2").Value = 0.07819433676692768
$ws.Range("P2").Value = 0.07819433676692769
$ws.Range("Q2").Value = 86.96598273688656
$ws.Range("R2").Value = 782.693844631979
$ws.Range("S2").Value = 0.01867905149923337
$ws.Range("T2").Value = 0.01867905149923337

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 14.80997966666667
$ws.Range("H3").Value = 44.429939
$ws.Range("I3").Value = 0.2388798507865045
$ws.Range("J3").Value = 0.2388798507865045
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 59.62659933333333
$ws.Range("N3").Value = 178.879798
$ws.Range("O3").Value = 0.7939998031155241
$ws.Range("P3").Value = 0.7939998031155242
$ws.Range("Q3").Value = 883.068723719147
$ws.Range("R3").Value = 7947.618513472323
$ws.Range("S3").Value = 0.1896705544927504
$ws.Range("T3").Value = 0.1896705544927504

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 14.80997966666667
$ws.Range("H4").Value = 44.429939
$ws.Range("I4").Value = 0.2388798507865045
$ws.Range("J4").Value = 0.2388798507865045
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.226786
$ws.Range("N4").Value = 0.680358
$ws.Range("O4").Value = 0.003019928041555994
$ws.Range("P4").Value = 0.003019928041555994
$ws.Range("Q4").Value = 3.358696048684668
$ws.Range("R4").Value = 30.228264438162
$ws.Range("S4").Value = 0.0007213999599528767
$ws.Range("T4").Value = 0.0007213999599528768

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 14.80997966666667
$ws.Range("H5").Value = 44.429939
$ws.Range("I5").Value = 0.2388798507865045
$ws.Range("J5").Value = 0.2388798507865045
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 9.370985666666668
$ws.Range("N5").Value = 28.112957
$ws.Range("O5").Value = 0.1247859320759922
$ws.Range("P5").Value = 0.1247859320759922
$ws.Range("Q5").Value = 138.7841071799581
$ws.Range("R5").Value = 1249.056964619623
$ws.Range("S5").Value = 0.02980884483456789
$ws.Range("T5").Value = 0.0298088448345679

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 33.42883866666667
$ws.Range("H6").Value = 100.286516
$ws.Range("I6").Value = 0.5391956081231261
$ws.Range("J6").Value = 0.5391956081231262
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 5.872120333333332
$ws.Range("N6").Value = 17.616361
$ws.Range("O6").Value = 0.07819433676692768
$ws.Range("P6").Value = 0.07819433676692769
$ws.Range("Q6").Value = 196.2981632542529
$ws.Range("R6").Value = 1766.683469288276
$ws.Range("S6").Value = 0.04216204296482808
$ws.Range("T6").Value = 0.0421620429648281

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 33.42883866666667
$ws.Range("H7").Value = 100.286516
$ws.Range("I7").Value = 0.5391956081231261
$ws.Range("J7").Value = 0.5391956081231262
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 59.62659933333333
$ws.Range("N7").Value = 178.879798
$ws.Range("O7").Value = 0.7939998031155241
$ws.Range("P7").Value = 0.7939998031155242
$ws.Range("Q7").Value = 1993.247969355974
$ws.Range("R7").Value = 17939.23172420377
$ws.Range("S7").Value = 0.4281212066905174
$ws.Range("T7").Value = 0.4281212066905175

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 33.42883866666667
$ws.Range("H8").Value = 100.286516
$ws.Range("I8").Value = 0.5391956081231261
$ws.Range("J8").Value = 0.5391956081231262
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.226786
$ws.Range("N8").Value = 0.680358
$ws.Range("O8").Value = 0.003019928041555994
$ws.Range("P8").Value = 0.003019928041555994
$ws.Range("Q8").Value = 7.581192605858668
$ws.Range("R8").Value = 68.230733452728
$ws.Range("S8").Value = 0.001628331936854865
$ws.Range("T8").Value = 0.001628331936854866

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 33.42883866666667
$ws.Range("H9").Value = 100.286516
$ws.Range("I9").Value = 0.5391956081231261
$ws.Range("J9").Value = 0.5391956081231262
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 9.370985666666668
$ws.Range("N9").Value = 28.112957
$ws.Range("O9").Value = 0.1247859320759922
$ws.Range("P9").Value = 0.1247859320759922
$ws.Range("Q9").Value = 313.2611679986459
$ws.Range("R9").Value = 2819.350511987812
$ws.Range("S9").Value = 0.0672840265309257
$ws.Range("T9").Value = 0.06728402653092572

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 2.509909333333333
$ws.Range("H10").Value = 7.529728
$ws.Range("I10").Value = 0.04048396962919451
$ws.Range("J10").Value = 0.04048396962919452
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 5.872120333333332
$ws.Range("N10").Value = 17.616361
$ws.Range("O10").Value = 0.07819433676692768
$ws.Range("P10").Value = 0.07819433676692769
$ws.Range("Q10").Value = 14.73848963108978
$ws.Range("R10").Value = 132.646406679808
$ws.Range("S10").Value = 0.003165617154847308
$ws.Range("T10").Value = 0.003165617154847309

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 2.509909333333333
$ws.Range("H11").Value = 7.529728
$ws.Range("I11").Value = 0.04048396962919451
$ws.Range("J11").Value = 0.04048396962919452
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 59.62659933333333
$ws.Range("N11").Value = 178.879798
$ws.Range("O11").Value = 0.7939998031155241
$ws.Range("P11").Value = 0.7939998031155242
$ws.Range("Q11").Value = 149.6573581816604
$ws.Range("R11").Value = 1346.916223634944
$ws.Range("S11").Value = 0.0321442639149153
$ws.Range("T11").Value = 0.03214426391491531

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 2.509909333333333
$ws.Range("H12").Value = 7.529728
$ws.Range("I12").Value = 0.04048396962919451
$ws.Range("J12").Value = 0.04048396962919452
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 0.226786
$ws.Range("N12").Value = 0.680358
$ws.Range("O12").Value = 0.003019928041555994
$ws.Range("P12").Value = 0.003019928041555994
$ws.Range("Q12").Value = 0.5692122980693334
$ws.Range("R12").Value = 5.122910682624
$ws.Range("S12").Value = 0.0001222586751167057
$ws.Range("T12").Value = 0.0001222586751167057

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 2.509909333333333
$ws.Range("H13").Value = 7.529728
$ws.Range("I13").Value = 0.04048396962919451
$ws.Range("J13").Value = 0.04048396962919452
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 9.370985666666668
$ws.Range("N13").Value = 28.112957
$ws.Range("O13").Value = 0.1247859320759922
$ws.Range("P13").Value = 0.1247859320759922
$ws.Range("Q13").Value = 23.52032438729956
$ws.Range("R13").Value = 211.682919485696
$ws.Range("S13").Value = 0.005051829884315196
$ws.Range("T13").Value = 0.005051829884315197

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 11.24888166666667
$ws.Range("H14").Value = 33.746645
$ws.Range("I14").Value = 0.1814405714611748
$ws.Range("J14").Value = 0.1814405714611748
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 5.872120333333332
$ws.Range("N14").Value = 17.616361
$ws.Range("O14").Value = 0.07819433676692768
$ws.Range("P14").Value = 0.07819433676692769
$ws.Range("Q14").Value = 66.05478676209388
$ws.Range("R14").Value = 594.4930808588449
$ws.Range("S14").Value = 0.01418762514801891
$ws.Range("T14").Value = 0.01418762514801892

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 11.24888166666667
$ws.Range("H15").Value = 33.746645
$ws.Range("I15").Value = 0.1814405714611748
$ws.Range("J15").Value = 0.1814405714611748
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 59.62659933333333
$ws.Range("N15").Value = 178.879798
$ws.Range("O15").Value = 0.7939998031155241
$ws.Range("P15").Value = 0.7939998031155242
$ws.Range("Q15").Value = 670.7325600864123
$ws.Range("R15").Value = 6036.59304077771
$ws.Range("S15").Value = 0.144063778017341
$ws.Range("T15").Value = 0.144063778017341

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 11.24888166666667
$ws.Range("H16").Value = 33.746645
$ws.Range("I16").Value = 0.1814405714611748
$ws.Range("J16").Value = 0.1814405714611748
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 0.226786
$ws.Range("N16").Value = 0.680358
$ws.Range("O16").Value = 0.003019928041555994
$ws.Range("P16").Value = 0.003019928041555994
$ws.Range("Q16").Value = 2.551088877656667
$ws.Range("R16").Value = 22.95979989891
$ws.Range("S16").Value = 0.0005479374696315461
$ws.Range("T16").Value = 0.0005479374696315462

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 11.24888166666667
$ws.Range("H17").Value = 33.746645
$ws.Range("I17").Value = 0.1814405714611748
$ws.Range("J17").Value = 0.1814405714611748
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 9.370985666666668
$ws.Range("N17").Value = 28.112957
$ws.Range("O17").Value = 0.1247859320759922
$ws.Range("P17").Value = 0.1247859320759922
$ws.Range("Q17").Value = 105.4131088643628
$ws.Range("R17").Value = 948.7179797792651
$ws.Range("S17").Value = 0.02264123082618336
$ws.Range("T17").Value = 0.02264123082618337
